$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (HP) ---
# Davi HP base (text-typed numeric string)
$ws.Range("D4").Value = "'939"
# Davi HP current (numeric)
$ws.Range("E4").Value = 932
# Lucca HP current (numeric)
$ws.Range("I4").Value = 61
# Magus HP current: was numeric 60, becomes text "110"
$ws.Range("AC4").Value = "'110"

# --- Row 5 (MP) ---
# Davi MP base + current (both text-typed numeric strings)
$ws.Range("D5").Value = "'90"
$ws.Range("E5").Value = "'90"

# --- Row 6 (Atk) ---
# Davi Atk base (numeric)
$ws.Range("D6").Value = 99

# --- Row 7 (Def) ---
# Davi / Lucca / Magus Def current (numeric)
$ws.Range("E7").Value = 0.44
$ws.Range("I7").Value = 0.3
$ws.Range("AC7").Value = 1.04

# --- Row 8 (Level) ---
$ws.Range("D8").Value = "'50"

# --- Row 9 (PWR) ---
$ws.Range("D9").Value = "'73"

# --- Row 11 (HIT) ---
$ws.Range("D11").Value = "'22"

# --- Row 12 (EV) ---
$ws.Range("D12").Value = "'26"

# --- Row 13 (STM) ---
$ws.Range("D13").Value = "'88"

# --- Row 14 (MAG) ---
$ws.Range("D14").Value = "'25"

# --- Row 15 (MDEF) ---
$ws.Range("D15").Value = "'80"
$ws.Range("E15").Value = 0.44
$ws.Range("I15").Value = 0.3
$ws.Range("AC15").Value = 1.04
